$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Copy the date/file/lines formatting from the first block (B7:D7) into the
# new K. Bicz (K7:M7) and M. Zygar (N7:P7) blocks before writing the values.
$ws.Range("B7:D7").Copy()
$ws.Range("K7:M7").PasteSpecial(-4122)
$ws.Range("N7:P7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New entry for "K. Bicz" (columns K:M) - row 7
$ws.Range("K7").Value = 45752
$ws.Range("L7").Value = "Diagram klas"
$ws.Range("M7").Value = 50

# New entry for "M. Zygar" (columns N:P) - row 7
$ws.Range("N7").Value = 45752
$ws.Range("O7").Value = "Diagram klas"
$ws.Range("P7").Value = 50

# Update active selection to match the saved workbook state
$null = $ws.Range("N9").Select()
